$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a text value into a cell while keeping it text-typed
# (forces NumberFormat to Text so Excel does not auto-convert numeric-looking
# strings like "6.80" or "1.00" into numbers, then resets the style back to
# Normal so no stray formatting/style index is introduced).
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") '64.465.25'
$ws.Range("E2").Value = '  -2.79%  '

# Row 3
Set-TextValue $ws.Range("D3") '3.177.54'
$ws.Range("E3").Value = '  -4.39%  '

# Row 5
Set-TextValue $ws.Range("D5") '571.17'
$ws.Range("E5").Value = '  -2.77%  '

# Row 6
Set-TextValue $ws.Range("D6") '168.98'
$ws.Range("E6").Value = '  -8.14%  '

# Row 7
Set-TextValue $ws.Range("D7") '0.607'
$ws.Range("E7").Value = '  -6.66%  '

# Row 8
$ws.Range("E8").Value = '  -0.10%  '

# Row 9
Set-TextValue $ws.Range("D9") '3.186.13'
$ws.Range("E9").Value = '  -4.07%  '

# Row 10
$ws.Range("E10").Value = '  -4.09%  '

# Row 11
Set-TextValue $ws.Range("D11") '6.80'
$ws.Range("E11").Value = '  -0.32%  '

# Row 12
$ws.Range("E12").Value = '  -3.75%  '

# Row 13
Set-TextValue $ws.Range("D13") '3.737.72'
$ws.Range("E13").Value = '  -4.15%  '

# Row 14
$ws.Range("E14").Value = '  -2.20%  '

# Row 15
Set-TextValue $ws.Range("D15") '64.512.17'
$ws.Range("E15").Value = '  -2.73%  '

# Row 16
Set-TextValue $ws.Range("D16") '25.40'
$ws.Range("E16").Value = '  -3.06%  '

# Row 17
$ws.Range("E17").Value = '  -2.93%  '

# Row 18
Set-TextValue $ws.Range("D18") '3.188.35'
$ws.Range("E18").Value = '  -4.23%  '

# Row 19
Set-TextValue $ws.Range("D19") '416.63'
$ws.Range("E19").Value = '  -2.22%  '

# Row 20
Set-TextValue $ws.Range("D20") '12.94'
$ws.Range("E20").Value = '  -2.20%  '

# Row 21
$ws.Range("E21").Value = '  -3.46%  '

# Row 22
Set-TextValue $ws.Range("D22") '7.11'
$ws.Range("E22").Value = '  -4.00%  '

# Row 23
Set-TextValue $ws.Range("D23") '0.999'
$ws.Range("E23").Value = '  -0.21%  '

# Row 24
Set-TextValue $ws.Range("D24") '70.50'
$ws.Range("E24").Value = '  -1.96%  '

# Row 25
Set-TextValue $ws.Range("D25") '5.68'
$ws.Range("E25").Value = '  -0.15%  '

# Row 26
$ws.Range("E26").Value = '  +0.69%  '

# Row 27
$ws.Range("E27").Value = '  -5.14%  '

# Row 28
$ws.Range("E28").Value = '  -7.05%  '

# Row 29
Set-TextValue $ws.Range("D29") '8.78'
$ws.Range("E29").Value = '  -1.86%  '

# Row 30
Set-TextValue $ws.Range("D30") '1.00'
$ws.Range("E30").Value = '  +0.06%  '

# Row 31
$ws.Range("E31").Value = '  -3.93%  '

# Row 32
Set-TextValue $ws.Range("D32") '21.76'
$ws.Range("E32").Value = '  -2.95%  '

# Row 33
$ws.Range("E33").Value = '  -0.12%  '

# Row 34
Set-TextValue $ws.Range("D34") '5.11'
$ws.Range("E34").Value = '  -1.70%  '

# Row 35
$ws.Range("E35").Value = '  -3.84%  '

# Row 36
$ws.Range("E36").Value = '  -3.87%  '

# Row 37
Set-TextValue $ws.Range("D37") '158.45'
$ws.Range("E37").Value = '  -0.93%  '

# Row 38
Set-TextValue $ws.Range("D38") '1.36'
$ws.Range("E38").Value = '  -5.55%  '

# Row 39
Set-TextValue $ws.Range("D39") '2.735.06'
$ws.Range("E39").Value = '  -5.21%  '

# Row 40
$ws.Range("E40").Value = '  -5.46%  '

# Row 41
Set-TextValue $ws.Range("D41") '24.43'
$ws.Range("E41").Value = '  -7.64%  '

# Row 42
Set-TextValue $ws.Range("D42") '4.20'
$ws.Range("E42").Value = '  -2.79%  '

# Row 43
Set-TextValue $ws.Range("D43") '39.20'
$ws.Range("E43").Value = '  -2.21%  '

# Row 44
Set-TextValue $ws.Range("D44") '0.716'
$ws.Range("E44").Value = '  -6.41%  '

# Row 45
$ws.Range("B45").Value = 'Hedera'
$ws.Range("C45").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range("D45") '0.0623'
$ws.Range("E45").Value = '  -6.34%  '

# Row 46
$ws.Range("B46").Value = 'RenderToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range("D46") '5.63'
$ws.Range("E46").Value = '  -5.02%  '

# Row 47
$ws.Range("E47").Value = '  -2.88%  '

# Row 48
Set-TextValue $ws.Range("D48") '21.70'
$ws.Range("E48").Value = '  -6.85%  '

# Row 49
Set-TextValue $ws.Range("D49") '293.63'
$ws.Range("E49").Value = '  -6.59%  '

# Row 50
Set-TextValue $ws.Range("D50") '2.02'
$ws.Range("E50").Value = '  -12.59%  '

# Row 51
$ws.Range("E51").Value = '  -0.19%  '
